# Weekly cryptos-list refresh: update each coin's Price (col D) and
# 1h Volume change (col E) to the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.178.40"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.434.16"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'316.26"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "'88.73"
$ws.Range("E6").Value = "  -4.39%  "
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D10").Value = "'32.13"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").Value = "'0.0832"
$ws.Range("E11").Value = "  -3.90%  "
$ws.Range("E12").Value = "  -2.72%  "
$ws.Range("D13").Value = "2.808.35"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").Value = "'15.67"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "2.431.43"
$ws.Range("E16").Value = "  -2.34%  "
$ws.Range("D17").Value = "'0.772"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "41.115.01"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").Value = "0.0₃0922"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("D20").Value = "'6.23"
$ws.Range("E20").Value = "  -4.01%  "
$ws.Range("D21").Value = "'71.88"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "'11.02"
$ws.Range("E22").Value = "  -4.22%  "
$ws.Range("D23").Value = "'235.45"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").Value = "'2.68"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").Value = "'23.98"
$ws.Range("E27").Value = "  -3.61%  "
$ws.Range("D28").Value = "'2.21"
$ws.Range("E28").Value = "  -3.23%  "
$ws.Range("D29").Value = "'9.55"
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("D30").Value = "'34.72"
$ws.Range("E30").Value = "  -4.73%  "
$ws.Range("D31").Value = "'156.83"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  -4.89%  "
$ws.Range("D34").Value = "'2.51"
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("D35").Value = "'0.0743"
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("D36").Value = "'2.92"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "'16.59"
$ws.Range("E37").Value = "  -4.92%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").Value = "'1.77"
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("E40").Value = "  -3.40%  "
$ws.Range("E41").Value = "  -3.74%  "
$ws.Range("E42").Value = "  -6.95%  "
$ws.Range("D43").Value = "1.983.28"
$ws.Range("E44").Value = "  -3.54%  "
$ws.Range("D45").Value = "'18.28"
$ws.Range("E45").Value = "  -5.86%  "
$ws.Range("E46").Value = "  -5.24%  "
$ws.Range("D47").Value = "'9.48"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").Value = "2.670.35"
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("D49").Value = "'95.15"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("D50").Value = "'73.03"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").Value = "'51.35"
$ws.Range("E51").Value = "  -2.45%  "

# The apostrophe-forced cells pick up Excel's quote-prefix marker; clear
# formatting across the whole Price column (contiguous range, so every
# area is actually covered) so they fall back to the original, unstyled
# appearance. Values already written are unaffected by ClearFormats.
$ws.Range("D2:D51").ClearFormats()

